# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after "总计" (and before "2022-Q2"),
#    populated with the quarterly fund-holding detail.
# 2. Update the "总计" (totals) sheet: insert a new row right under the header
#    for "2022-Q3" and shift the previously existing rows down, renumbering the
#    leading index column.

function Set-TextCell($range, [string]$value) {
    # Force the cell to be stored as text even when the value looks numeric
    # (fund codes / percentages kept as strings in the source data), then
    # drop the "quote prefix" number-format side effect so the cell keeps the
    # plain/default style.
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")

# --- Part 1: add the "2022-Q3" worksheet -----------------------------------

$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# Source sheet to copy header / index-column formatting from.
$fmtSrc = $wb.Worksheets.Item("2022-Q2")

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
$col = 2
foreach ($h in $headers) {
    $q3.Cells.Item(1, $col).Value = $h
    $fmtSrc.Cells.Item(1, $col).Copy()
    $q3.Cells.Item(1, $col).PasteSpecial(-4122)
    $col++
}

$rows = @(
    @(0, "159869", "华夏中证动漫游戏ETF", "6.35", "99.31", "5.55", "0.3524", 7),
    @(1, "516010", "国泰中证动漫游戏ETF", "3.78", "97.86", "5.41", "0.2045", 7),
    @(2, "005585", "银河文体娱乐主题灵活配置混合A", "3.01", "90.28", "4.98", "0.1499", 5),
    @(3, "001628", "招商体育文化休闲股票A", "2.23", "92.42", "4.90", "0.1093", 8),
    @(4, "004809", "新疆前海联合润丰灵活配置混合A", "1.30", "87.05", "5.32", "0.0692", 6),
    @(5, "516770", "华泰柏瑞中证动漫游戏ETF", "0.99", "96.39", "5.50", "0.0544", 7),
    @(6, "161036", "富国中证娱乐主题指数增强（LOF）A", "0.74", "93.50", "3.08", "0.0228", 8),
    @(7, "517500", "国泰中证沪港深动漫游戏ETF", "0.53", "92.78", "3.88", "0.0206", 8),
    @(8, "015667", "银河文体娱乐主题灵活配置混合C", "0.41", "90.28", "4.98", "0.0204", 5),
    @(9, "004890", "中邮健康文娱灵活配置混合", "0.41", "86.15", "3.52", "0.0144", 10),
    @(10, "015395", "招商体育文化休闲股票C", "0.25", "92.42", "4.90", "0.0122", 8),
    @(11, "005167", "嘉实润泽量化一年定期开放混合", "0.55", "24.55", "0.66", "0.0036", 2),
    @(12, "014256", "富国中证娱乐主题指数增强（LOF）C", "0.07", "93.50", "3.08", "0.0022", 8),
    @(13, "005935", "新疆前海联合润丰灵活配置混合C", "0.01", "87.05", "5.32", "0.0005", 6)
)

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $fmtSrc.Cells.Item(2, 1).Copy()
    $q3.Cells.Item($r, 1).PasteSpecial(-4122)

    Set-TextCell $q3.Cells.Item($r, 2) $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    Set-TextCell $q3.Cells.Item($r, 4) $row[3]
    Set-TextCell $q3.Cells.Item($r, 5) $row[4]
    Set-TextCell $q3.Cells.Item($r, 6) $row[5]
    Set-TextCell $q3.Cells.Item($r, 7) $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# --- Part 2: update the "总计" sheet ----------------------------------------

$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("B2:D2").Style = "Normal"
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 14
$totalSheet.Range("D2").Value = 1.04

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
